$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Aerial Systems - Applications I"
$ws.Range("D6").Value = "Aerial Systems - Applications II"
$ws.Range("D7").Value = "Aerial Systems - Applications III"

$ws.Columns("D:D").ColumnWidth = 28.5

$ws.Range("E6").Select() | Out-Null
